$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3e"
$ws.Range("C2").Value = "Plxnd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.032802
$ws.Range("H2").Value = 0.098406
$ws.Range("I2").Value = 0.02084482040973249
$ws.Range("J2").Value = 0.02084482040973249
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 139.2986196666667
$ws.Range("N2").Value = 417.895859
$ws.Range("O2").Value = 0.6137320738580456
$ws.Range("P2").Value = 0.6137320738580456
$ws.Range("Q2").Value = 4.569273322306
$ws.Range("R2").Value = 41.12345990075399
$ws.Range("S2").Value = 0.01279313485926364
$ws.Range("T2").Value = 0.01279313485926364

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3e"
$ws.Range("C3").Value = "Plxnd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.032802
$ws.Range("H3").Value = 0.098406
$ws.Range("I3").Value = 0.02084482040973249
$ws.Range("J3").Value = 0.02084482040973249
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 22.17197066666667
$ws.Range("N3").Value = 66.515912
$ws.Range("O3").Value = 0.09768689432339951
$ws.Range("P3").Value = 0.09768689432339951
$ws.Range("Q3").Value = 0.7272849818079999
$ws.Range("R3").Value = 6.545564836272
$ws.Range("S3").Value = 0.002036265768555779
$ws.Range("T3").Value = 0.00203626576855578

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3e"
$ws.Range("C4").Value = "Plxnd1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.032802
$ws.Range("H4").Value = 0.098406
$ws.Range("I4").Value = 0.02084482040973249
$ws.Range("J4").Value = 0.02084482040973249
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 35.78898466666666
$ws.Range("N4").Value = 107.366954
$ws.Range("O4").Value = 0.1576817331952585
$ws.Range("P4").Value = 0.1576817331952585
$ws.Range("Q4").Value = 1.173950275036
$ws.Range("R4").Value = 10.565552475324
$ws.Range("S4").Value = 0.003286847410350518
$ws.Range("T4").Value = 0.003286847410350519

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema3e"
$ws.Range("C5").Value = "Plxnd1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.032802
$ws.Range("H5").Value = 0.098406
$ws.Range("I5").Value = 0.02084482040973249
$ws.Range("J5").Value = 0.02084482040973249
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 29.710182
$ws.Range("N5").Value = 89.13054600000001
$ws.Range("O5").Value = 0.1308992986232963
$ws.Range("P5").Value = 0.1308992986232963
$ws.Range("Q5").Value = 0.9745533899640001
$ws.Range("R5").Value = 8.770980509676
$ws.Range("S5").Value = 0.002728572371562555
$ws.Range("T5").Value = 0.002728572371562556

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3e"
$ws.Range("C6").Value = "Plxnd1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1716996666666667
$ws.Range("H6").Value = 0.5150990000000001
$ws.Range("I6").Value = 0.1091106858142064
$ws.Range("J6").Value = 0.1091106858142065
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 139.2986196666667
$ws.Range("N6").Value = 417.895859
$ws.Range("O6").Value = 0.6137320738580456
$ws.Range("P6").Value = 0.6137320738580456
$ws.Range("Q6").Value = 23.91752656389345
$ws.Range("R6").Value = 215.257739075041
$ws.Range("S6").Value = 0.06696472748482656
$ws.Range("T6").Value = 0.06696472748482657

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3e"
$ws.Range("C7").Value = "Plxnd1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1716996666666667
$ws.Range("H7").Value = 0.5150990000000001
$ws.Range("I7").Value = 0.1091106858142064
$ws.Range("J7").Value = 0.1091106858142065
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.17197066666667
$ws.Range("N7").Value = 66.515912
$ws.Range("O7").Value = 0.09768689432339951
$ws.Range("P7").Value = 0.09768689432339951
$ws.Range("Q7").Value = 3.806919972809778
$ws.Range("R7").Value = 34.26227975528801
$ws.Range("S7").Value = 0.01065868403468603
$ws.Range("T7").Value = 0.01065868403468603

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema3e"
$ws.Range("C8").Value = "Plxnd1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1716996666666667
$ws.Range("H8").Value = 0.5150990000000001
$ws.Range("I8").Value = 0.1091106858142064
$ws.Range("J8").Value = 0.1091106858142065
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 35.78898466666666
$ws.Range("N8").Value = 107.366954
$ws.Range("O8").Value = 0.1576817331952585
$ws.Range("P8").Value = 0.1576817331952585
$ws.Range("Q8").Value = 6.144956737605112
$ws.Range("R8").Value = 55.30461063844601
$ws.Range("S8").Value = 0.01720476204930738
$ws.Range("T8").Value = 0.01720476204930739

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema3e"
$ws.Range("C9").Value = "Plxnd1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1716996666666667
$ws.Range("H9").Value = 0.5150990000000001
$ws.Range("I9").Value = 0.1091106858142064
$ws.Range("J9").Value = 0.1091106858142065
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 29.710182
$ws.Range("N9").Value = 89.13054600000001
$ws.Range("O9").Value = 0.1308992986232963
$ws.Range("P9").Value = 0.1308992986232963
$ws.Range("Q9").Value = 5.101228346006002
$ws.Range("R9").Value = 45.91105511405402
$ws.Range("S9").Value = 0.01428251224538647
$ws.Range("T9").Value = 0.01428251224538647

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema3e"
$ws.Range("C10").Value = "Plxnd1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.365476
$ws.Range("H10").Value = 4.096428
$ws.Range("I10").Value = 0.8677245897750103
$ws.Range("J10").Value = 0.8677245897750105
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 139.2986196666667
$ws.Range("N10").Value = 417.895859
$ws.Range("O10").Value = 0.6137320738580456
$ws.Range("P10").Value = 0.6137320738580456
$ws.Range("Q10").Value = 190.2089219879614
$ws.Range("R10").Value = 1711.880297891652
$ws.Range("S10").Value = 0.532550412020239
$ws.Range("T10").Value = 0.5325504120202391

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema3e"
$ws.Range("C11").Value = "Plxnd1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.365476
$ws.Range("H11").Value = 4.096428
$ws.Range("I11").Value = 0.8677245897750103
$ws.Range("J11").Value = 0.8677245897750105
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 22.17197066666667
$ws.Range("N11").Value = 66.515912
$ws.Range("O11").Value = 0.09768689432339951
$ws.Range("P11").Value = 0.09768689432339951
$ws.Range("Q11").Value = 30.27529381803734
$ws.Range("R11").Value = 272.477644362336
$ws.Range("S11").Value = 0.08476532030316662
$ws.Range("T11").Value = 0.08476532030316665

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema3e"
$ws.Range("C12").Value = "Plxnd1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.365476
$ws.Range("H12").Value = 4.096428
$ws.Range("I12").Value = 0.8677245897750103
$ws.Range("J12").Value = 0.8677245897750105
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 35.78898466666666
$ws.Range("N12").Value = 107.366954
$ws.Range("O12").Value = 0.1576817331952585
$ws.Range("P12").Value = 0.1576817331952585
$ws.Range("Q12").Value = 48.86899962670133
$ws.Range("R12").Value = 439.820996640312
$ws.Range("S12").Value = 0.1368243172518683
$ws.Range("T12").Value = 0.1368243172518684

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema3e"
$ws.Range("C13").Value = "Plxnd1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.365476
$ws.Range("H13").Value = 4.096428
$ws.Range("I13").Value = 0.8677245897750103
$ws.Range("J13").Value = 0.8677245897750105
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 29.710182
$ws.Range("N13").Value = 89.13054600000001
$ws.Range("O13").Value = 0.1308992986232963
$ws.Range("P13").Value = 0.1308992986232963
$ws.Range("Q13").Value = 40.56854047663201
$ws.Range("R13").Value = 365.1168642896881
$ws.Range("S13").Value = 0.1135845401997364
$ws.Range("T13").Value = 0.1135845401997364

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Sema3e"
$ws.Range("C14").Value = "Plxnd1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.003650666666666667
$ws.Range("H14").Value = 0.010952
$ws.Range("I14").Value = 0.00231990400105065
$ws.Range("J14").Value = 0.00231990400105065
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 139.2986196666667
$ws.Range("N14").Value = 417.895859
$ws.Range("O14").Value = 0.6137320738580456
$ws.Range("P14").Value = 0.6137320738580456
$ws.Range("Q14").Value = 0.5085328275297778
$ws.Range("R14").Value = 4.576795447767999
$ws.Range("S14").Value = 0.001423799493716393
$ws.Range("T14").Value = 0.001423799493716393

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Sema3e"
$ws.Range("C15").Value = "Plxnd1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.003650666666666667
$ws.Range("H15").Value = 0.010952
$ws.Range("I15").Value = 0.00231990400105065
$ws.Range("J15").Value = 0.00231990400105065
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 22.17197066666667
$ws.Range("N15").Value = 66.515912
$ws.Range("O15").Value = 0.09768689432339951
$ws.Range("P15").Value = 0.09768689432339951
$ws.Range("Q15").Value = 0.08094247424711111
$ws.Range("R15").Value = 0.728482268224
$ws.Range("S15").Value = 0.0002266242169910665
$ws.Range("T15").Value = 0.0002266242169910666

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Sema3e"
$ws.Range("C16").Value = "Plxnd1"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.003650666666666667
$ws.Range("H16").Value = 0.010952
$ws.Range("I16").Value = 0.00231990400105065
$ws.Range("J16").Value = 0.00231990400105065
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 35.78898466666666
$ws.Range("N16").Value = 107.366954
$ws.Range("O16").Value = 0.1576817331952585
$ws.Range("P16").Value = 0.1576817331952585
$ws.Range("Q16").Value = 0.1306536533564444
$ws.Range("R16").Value = 1.175882880208
$ws.Range("S16").Value = 0.0003658064837322813
$ws.Range("T16").Value = 0.0003658064837322814

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema3e"
$ws.Range("C17").Value = "Plxnd1"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.003650666666666667
$ws.Range("H17").Value = 0.010952
$ws.Range("I17").Value = 0.00231990400105065
$ws.Range("J17").Value = 0.00231990400105065
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 29.710182
$ws.Range("N17").Value = 89.13054600000001
$ws.Range("O17").Value = 0.1308992986232963
$ws.Range("P17").Value = 0.1308992986232963
$ws.Range("Q17").Value = 0.108461971088
$ws.Range("R17").Value = 0.976157739792
$ws.Range("S17").Value = 0.000303673806610909
$ws.Range("T17").Value = 0.0003036738066109091

